$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.802.14'
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").Value = '1.800.33'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''308.75'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '''0.4652'
$ws.Range("E7").Value = '  +4.04%  '
$ws.Range("D8").Value = '''0.3674'
$ws.Range("E8").Value = '  -2.10%  '
$ws.Range("D9").Value = '''0.07358'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = '''0.8655'
$ws.Range("E10").Value = '  -2.35%  '
$ws.Range("D11").Value = '''20.35'
$ws.Range("E11").Value = '  -3.18%  '
$ws.Range("D12").Value = '1.857.20'
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("D14").Value = '''6.519'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").Value = '''0.07038'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").Value = '''91.21'
$ws.Range("E16").Value = '  -2.81%  '
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '''0.000008704'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = '''1.001'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '''14.62'
$ws.Range("E20").Value = '  -3.53%  '
$ws.Range("D21").Value = '26.826.26'
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("D22").Value = '''5.281'
$ws.Range("E22").Value = '  -2.52%  '
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("D24").Value = '2.132.36'
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").Value = '''1.901'
$ws.Range("E25").Value = '  -3.27%  '
$ws.Range("D26").Value = '''151.24'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '''18.32'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("D28").Value = '''2.121'
$ws.Range("E28").Value = '  -8.00%  '
$ws.Range("D29").Value = '''5.231'
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("D30").Value = '''115.75'
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("D31").Value = '''0.08903'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").Value = '''0.7551'
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").Value = '''1.147'
$ws.Range("E34").Value = '  -4.81%  '
$ws.Range("E35").Value = '  -3.94%  '
$ws.Range("D36").Value = '''1.001'
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '''1.106'
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("D38").Value = '''0.01949'
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("D39").Value = '''0.05239'
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").Value = '''2.928'
$ws.Range("E40").Value = '  +2.55%  '
$ws.Range("D41").Value = '''7.192'
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("D42").Value = '''0.5261'
$ws.Range("E42").Value = '  -1.78%  '
$ws.Range("D43").Value = '''2.326'
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("D44").Value = '''0.1657'
$ws.Range("E44").Value = '  -3.49%  '
$ws.Range("D45").Value = '''8.446'
$ws.Range("E45").Value = '  -2.57%  '
$ws.Range("D46").Value = '''0.4993'
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("D47").Value = '''10.28'
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''104.02'
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").Value = '''1.001'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '''1.661'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("D51").Value = '''0.06284'
$ws.Range("E51").Value = '  -1.85%  '
